$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force the value to be stored as text (avoids Excel's automatic
# number coercion for numeric-looking strings like "250.40" or "1.00"),
# then restore the cell's style to Normal so no stray number format sticks.
function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 (Bitcoin)
Set-TextValue "D2" "36.787.98"
Set-TextValue "E2" "  +4.16%  "

# Row 3 (Ethereum)
Set-TextValue "D3" "1.912.84"
Set-TextValue "E3" "  +1.66%  "

# Row 4 (TetherUSD)
Set-TextValue "E4" "  -0.01%  "

# Row 5 and Row 6: swapped coins (XRP <-> BNB)
Set-TextValue "B5" "BNB"
Set-TextValue "C5" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D5" "250.40"
Set-TextValue "E5" "  +1.68%  "

Set-TextValue "B6" "XRP"
Set-TextValue "C6" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D6" "0.704"
Set-TextValue "E6" "  +3.14%  "

# Row 7 (USDC)
Set-TextValue "D7" "1.00"
Set-TextValue "E7" "  -0.03%  "

# Row 8 (Solana)
Set-TextValue "D8" "46.51"
Set-TextValue "E8" "  +7.29%  "

# Row 9 (Cardano)
Set-TextValue "D9" "0.372"
Set-TextValue "E9" "  +4.63%  "

# Row 10 (OKB)
Set-TextValue "D10" "58.55"

# Row 11 (Dogecoin)
Set-TextValue "D11" "0.0764"
Set-TextValue "E11" "  +1.93%  "

# Row 12 (TRON)
Set-TextValue "D12" "0.0998"
Set-TextValue "E12" "  +2.11%  "

# Row 13 (Chainlink)
Set-TextValue "D13" "14.69"
Set-TextValue "E13" "  +8.63%  "

# Row 14 (Polygon)
Set-TextValue "D14" "0.815"
Set-TextValue "E14" "  +7.17%  "

# Row 15 (WrappedliquidstakedEther2.0)
Set-TextValue "E15" "  +1.73%  "

# Row 16 (Polkadot)
Set-TextValue "D16" "5.14"
Set-TextValue "E16" "  +4.09%  "

# Row 17 (WrappedEther)
Set-TextValue "D17" "1.923.86"
Set-TextValue "E17" "  +2.58%  "

# Row 18 (WrappedBTC)
Set-TextValue "D18" "36.771.95"
Set-TextValue "E18" "  +4.06%  "

# Row 19 (Litecoin)
Set-TextValue "D19" "74.95"

# Row 20 (ShibaInu)
Set-TextValue "E20" "  +4.52%  "

# Row 21 (BitcoinCash)
Set-TextValue "D21" "251.20"
Set-TextValue "E21" "  +2.86%  "

# Row 22 (Avalanche)
Set-TextValue "D22" "13.41"
Set-TextValue "E22" "  +4.80%  "

# Row 23 (Uniswap)
Set-TextValue "D23" "5.18"
Set-TextValue "E23" "  -0.11%  "

# Row 24 (Toncoin)
Set-TextValue "E24" "  +0.71%  "

# Row 25 (Dai)
Set-TextValue "E25" "  +0.07%  "

# Row 26 (PancakeSwap)
Set-TextValue "D26" "2.19"
Set-TextValue "E26" "  +0.47%  "

# Row 27 (Monero)
Set-TextValue "D27" "168.13"
Set-TextValue "E27" "  +2.06%  "

# Row 28 (Cosmos)
Set-TextValue "E28" "  +1.91%  "

# Row 29 (EthereumClassic)
Set-TextValue "D29" "18.79"
Set-TextValue "E29" "  +2.68%  "

# Row 30 (Stellar)
Set-TextValue "E30" "  +1.82%  "

# Row 31 (Filecoin)
Set-TextValue "D31" "4.56"
Set-TextValue "E31" "  +6.04%  "

# Row 32 (Hedera)
Set-TextValue "D32" "0.0617"
Set-TextValue "E32" "  +3.96%  "

# Row 33 (InternetComputer)
Set-TextValue "E33" "  +3.41%  "

# Row 34 (Kaspa)
Set-TextValue "D34" "0.0896"
Set-TextValue "E34" "  +23.23%  "

# Row 35 (BinanceUSD)
Set-TextValue "E35" "  -0.02%  "

# Row 36 (WEMIXToken)
Set-TextValue "D36" "1.86"
Set-TextValue "E36" "  +0.95%  "

# Row 37 (TrustWalletToken)
Set-TextValue "E37" "  +4.94%  "

# Row 38 (ImmutableX)
Set-TextValue "D38" "0.873"
Set-TextValue "E38" "  +2.44%  "

# Row 39 (Gas)
Set-TextValue "D39" "17.94"
Set-TextValue "E39" "  +53.27%  "

# Row 40 (LidoDAOToken)
Set-TextValue "E40" "  +2.96%  "

# Row 41 (Aave)
Set-TextValue "D41" "105.11"
Set-TextValue "E41" "  +8.31%  "

# Row 42 (VeChain)
Set-TextValue "E42" "  +3.89%  "

# Row 43 (InjectiveProtocol)
Set-TextValue "D43" "17.64"
Set-TextValue "E43" "  +0.80%  "

# Row 44 (HuobiToken)
Set-TextValue "E44" "  +21.96%  "

# Row 45 (ARBITRUM)
Set-TextValue "E45" "  +2.80%  "

# Row 46 (Maker)
Set-TextValue "D46" "1.350.73"
Set-TextValue "E46" "  +3.14%  "

# Row 47 (RenderToken)
Set-TextValue "D47" "2.37"
Set-TextValue "E47" "  -1.56%  "

# Row 48 (Cronos)
Set-TextValue "E48" "  +1.65%  "

# Row 49 (MXToken)
Set-TextValue "E49" "  +2.52%  "

# Row 50 (FraxShare)
Set-TextValue "E50" "  +2.29%  "

# Row 51 (MultiversX)
Set-TextValue "E51" "  +2.78%  "
